# Scheduled-runner refresh: push updated market-price snapshots (and the
# profit figures derived from them) into the per-job sheets. Pure value
# overwrites - no formulas/formatting involved, matching how the source
# data is produced upstream.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 385607.2
$ws.Range("J17").Value = 385607.2
$ws.Range("L17").Value = 1156821.6
$ws.Range("N17").Value = -1157157.6

$ws.Range("H28").Value = 737.6316
$ws.Range("I28").Value = 886.8
$ws.Range("J28").Value = 178.25
$ws.Range("K28").Value = 886.8
$ws.Range("L28").Value = 178.25
$ws.Range("M28").Value = -401.8
$ws.Range("N28").Value = -1148.25

$ws.Range("H40").Value = 1879.1305
$ws.Range("I40").Value = 1605.4445
$ws.Range("K40").Value = 1605.4445
$ws.Range("M40").Value = -1430.4445

$ws.Range("H112").Value = 1333.4445
$ws.Range("I112").Value = 949.8570999999999
$ws.Range("J112").Value = 1467.7
$ws.Range("K112").Value = 2849.5713
$ws.Range("L112").Value = 4403.1
$ws.Range("M112").Value = -1741.5713
$ws.Range("N112").Value = -6619.1

$ws.Range("H116").Value = 4713.3076
$ws.Range("I116").Value = 4919.222
$ws.Range("J116").Value = 4250
$ws.Range("K116").Value = 4919.222
$ws.Range("L116").Value = 4250
$ws.Range("M116").Value = -1477.222
$ws.Range("N116").Value = -11134

$ws.Range("H129").Value = 1195609.9
$ws.Range("J129").Value = 2058921.1
$ws.Range("L129").Value = 6176763.300000001
$ws.Range("N129").Value = -6186763.300000001

$ws.Range("H138").Value = 2457.23
$ws.Range("I138").Value = 830.88
$ws.Range("J138").Value = 2999.3467
$ws.Range("K138").Value = 2492.64
$ws.Range("L138").Value = 8998.0401
$ws.Range("M138").Value = 2647.36
$ws.Range("N138").Value = -19278.0401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3719241
$ws.Range("I32").Value = 4183334.8
$ws.Range("J32").Value = 6491.4
$ws.Range("K32").Value = 4183334.8
$ws.Range("L32").Value = 6491.4
$ws.Range("M32").Value = -4183047.8
$ws.Range("N32").Value = -7065.4

$ws.Range("H45").Value = 1570.9048
$ws.Range("I45").Value = 894.93335
$ws.Range("J45").Value = 3260.8333
$ws.Range("K45").Value = 894.93335
$ws.Range("L45").Value = 3260.8333
$ws.Range("M45").Value = -517.93335
$ws.Range("N45").Value = -4014.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28659.578
$ws.Range("I134").Value = 2108.6667
$ws.Range("K134").Value = 6326.000100000001
$ws.Range("M134").Value = -3791.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 235
$ws.Range("I19").Value = 235
$ws.Range("K19").Value = 235
$ws.Range("M19").Value = -65

$ws.Range("H24").Value = 235
$ws.Range("I24").Value = 235
$ws.Range("K24").Value = 235
$ws.Range("M24").Value = -65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 46.285713
$ws.Range("I12").Value = 136.66667
$ws.Range("J12").Value = 21.636364
$ws.Range("K12").Value = 410.00001
$ws.Range("L12").Value = 64.909092
$ws.Range("M12").Value = -237.00001
$ws.Range("N12").Value = -410.909092

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null

$ws.Range("H80").Value = 5782.143
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 7605.5557
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 22816.6671
$ws.Range("M80").Value = -6564
$ws.Range("N80").Value = -24688.6671

$ws.Range("H83").Value = 5782.143
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 7605.5557
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 68450.0013
$ws.Range("M83").Value = -17820
$ws.Range("N83").Value = -77810.0013

$ws.Range("H113").Value = 657.625
$ws.Range("J113").Value = 523
$ws.Range("L113").Value = 1569
$ws.Range("N113").Value = -5909

$ws.Range("H131").Value = 1050.0769
$ws.Range("I131").Value = 1678.75
$ws.Range("J131").Value = 989.48193
$ws.Range("K131").Value = 5036.25
$ws.Range("L131").Value = 2968.44579
$ws.Range("M131").Value = 3.75
$ws.Range("N131").Value = -13048.44579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4871.7856
$ws.Range("I70").Value = 4875
$ws.Range("J70").Value = 4870.5
$ws.Range("K70").Value = 4875
$ws.Range("L70").Value = 4870.5
$ws.Range("M70").Value = -4605
$ws.Range("N70").Value = -5410.5

$ws.Range("H73").Value = 4871.7856
$ws.Range("I73").Value = 4875
$ws.Range("J73").Value = 4870.5
$ws.Range("K73").Value = 4875
$ws.Range("L73").Value = 4870.5
$ws.Range("M73").Value = -3939
$ws.Range("N73").Value = -6742.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1934.3077
$ws.Range("I22").Value = 2375.1667
$ws.Range("J22").Value = 1556.4286
$ws.Range("K22").Value = 2375.1667
$ws.Range("L22").Value = 1556.4286
$ws.Range("M22").Value = -2080.1667
$ws.Range("N22").Value = -2146.4286

$ws.Range("H27").Value = 1934.3077
$ws.Range("I27").Value = 2375.1667
$ws.Range("J27").Value = 1556.4286
$ws.Range("K27").Value = 2375.1667
$ws.Range("L27").Value = 1556.4286
$ws.Range("M27").Value = -2268.1667
$ws.Range("N27").Value = -1770.4286

$ws.Range("H46").Value = 1229.45
$ws.Range("J46").Value = 1102.3334
$ws.Range("L46").Value = 1102.3334
$ws.Range("N46").Value = -1478.3334

$ws.Range("H96").Value = 13999
$ws.Range("J96").Value = 13999
$ws.Range("L96").Value = 13999
$ws.Range("N96").Value = -19491

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -771

$ws.Range("H99").Value = 42000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 42000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 42000
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = -47990

$ws.Range("H122").Value = 1370.9656
$ws.Range("I122").Value = 1336.8462
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 4010.5386
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -1560.5386
$ws.Range("N122").Value = -9899.9998
